$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "gidapp" bet block, mirroring the existing row-7 block's layout
# (A: date/time label, B: the 3-digit bet, C: amount).
$ws.Range("A12").Value = "09 sun jun 2019 2"
$ws.Range("B12").Value = 713
$ws.Range("C12").Value = 10

$ws.Range("B13").Value = 313
$ws.Range("C13").Value = 10

# Excel leaves the active selection on the last-entered cell.
$ws.Range("C13").Select()
